# Add the new "Sheet1" worksheet at the end of the workbook (after "35K EC")
# and populate it with the "Adopted 35Ar" working data.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Sheet1"

# Row 1
$ws.Range("B1").Value = "35AR"
$ws.Range("C1").Value = "G"
$ws.Range("D1").Value = 1162
$ws.Range("E1").Value = 8
$ws.Range("F1").Value = 15
$ws.Range("G1").Value = 3

# Row 2
$ws.Range("B2").Value = "35AR"
$ws.Range("C2").Value = "G"
$ws.Range("D2").Value = 1756.3
$ws.Range("E2").Value = 14
$ws.Range("F2").Value = 17
$ws.Range("G2").Value = 9

# Row 4 (row 3 left blank)
$ws.Range("B4").Value = "35AR"
$ws.Range("C4").Value = "G"
$ws.Range("D4").Value = 1162
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 11
$ws.Range("G4").Value = 3

# Row 5
$ws.Range("B5").Value = "35AR"
$ws.Range("C5").Value = "G"
$ws.Range("D5").Value = 1756
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 27
$ws.Range("G5").Value = 4
$ws.Range("H5").Value = "?"

# Match the author's final selection on the new sheet
[void]$ws.Range("H5").Select()

# Match the page setup of the other sheets in the workbook (A4, portrait)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
